$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "SubsymptomTable"
$wb.Worksheets.Item("Sheet1").Name = "SubsymptomTable"

# Adjust view on SymptomTable sheet: zoom and frozen pane scroll position
$wsSymptom = $wb.Worksheets.Item("SymptomTable")
$wsSymptom.Activate()
$excel.ActiveWindow.Zoom = 101
$wsSymptom.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$excel.ActiveWindow.FreezePanes = $true
$wsSymptom.Range("A233").Select()

# Adjust selection on SubsymptomTable sheet
$wsSub = $wb.Worksheets.Item("SubsymptomTable")
$wsSub.Activate()
$wsSub.Range("F13").Select()
